# Auto-generated Excel COM-interop script to apply the diff
# (scheduled-runner market price refresh) to Halicarnassus_Profits sheets.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 781.25
$ws.Range("I6").Value = 40.625
$ws.Range("K6").Value = 121.875
$ws.Range("M6").Value = -9.875
$ws.Range("H12").Value = 119.4
$ws.Range("I12").Value = 73
$ws.Range("J12").Value = 189
$ws.Range("K12").Value = 73
$ws.Range("L12").Value = 189
$ws.Range("M12").Value = 97
$ws.Range("N12").Value = -529
$ws.Range("H40").Value = 5802.1177
$ws.Range("I40").Value = 4444
$ws.Range("J40").Value = 7742.2856
$ws.Range("K40").Value = 4444
$ws.Range("L40").Value = 7742.2856
$ws.Range("M40").Value = -4269
$ws.Range("N40").Value = -8092.2856
$ws.Range("H70").Value = 3827.9
$ws.Range("J70").Value = 4142.222
$ws.Range("L70").Value = 12426.666
$ws.Range("N70").Value = -12966.666
$ws.Range("H73").Value = 3827.9
$ws.Range("J73").Value = 4142.222
$ws.Range("L73").Value = 12426.666
$ws.Range("N73").Value = -14298.666
$ws.Range("H87").Value = 69661
$ws.Range("J87").Value = 69661
$ws.Range("L87").Value = 69661
$ws.Range("N87").Value = -72157
$ws.Range("H90").Value = 69661
$ws.Range("J90").Value = 69661
$ws.Range("L90").Value = 208983
$ws.Range("N90").Value = -221463
$ws.Range("H115").Value = 1333.6364
$ws.Range("I115").Value = 963.44446
$ws.Range("J115").Value = 2999.5
$ws.Range("K115").Value = 2890.33338
$ws.Range("L115").Value = 8998.5
$ws.Range("M115").Value = -1323.33338
$ws.Range("N115").Value = -12132.5
$ws.Range("H117").Value = 109995
$ws.Range("J117").Value = 109995
$ws.Range("L117").Value = 109995
$ws.Range("N117").Value = -119173
$ws.Range("H135").Value = 2470.4443
$ws.Range("I135").Value = 2483.1667
$ws.Range("J135").Value = 2445
$ws.Range("K135").Value = 22348.5003
$ws.Range("L135").Value = 22005
$ws.Range("M135").Value = -19813.5003
$ws.Range("N135").Value = -27075

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2048.3845
$ws.Range("I2").Value = 1572.9
$ws.Range("J2").Value = 3633.3333
$ws.Range("K2").Value = 1572.9
$ws.Range("L2").Value = 3633.3333
$ws.Range("M2").Value = -1459.9
$ws.Range("N2").Value = -3859.3333
$ws.Range("H62").Value = 50000
$ws.Range("J62").Value = 50000
$ws.Range("L62").Value = 50000
$ws.Range("N62").Value = -51248
$ws.Range("H63").Value = 4146.2
$ws.Range("I63").Value = 2475.75
$ws.Range("J63").Value = 5259.8335
$ws.Range("K63").Value = 2475.75
$ws.Range("L63").Value = 5259.8335
$ws.Range("M63").Value = -1789.75
$ws.Range("N63").Value = -6631.8335
$ws.Range("H65").Value = 50000
$ws.Range("J65").Value = 50000
$ws.Range("L65").Value = 150000
$ws.Range("N65").Value = -156240
$ws.Range("H66").Value = 4146.2
$ws.Range("I66").Value = 2475.75
$ws.Range("J66").Value = 5259.8335
$ws.Range("K66").Value = 12378.75
$ws.Range("L66").Value = 26299.1675
$ws.Range("M66").Value = -8946.75
$ws.Range("N66").Value = -33163.1675
$ws.Range("H110").Value = 3416.889
$ws.Range("I110").Value = 3002
$ws.Range("K110").Value = 3002
$ws.Range("M110").Value = -957
$ws.Range("H116").Value = 2048.3845
$ws.Range("I116").Value = 1572.9
$ws.Range("J116").Value = 3633.3333
$ws.Range("K116").Value = 1572.9
$ws.Range("L116").Value = 3633.3333
$ws.Range("M116").Value = 721.0999999999999
$ws.Range("N116").Value = -8221.3333
$ws.Range("H132").Value = 1991.1177
$ws.Range("I132").Value = 1282.5
$ws.Range("J132").Value = 5298
$ws.Range("K132").Value = 3847.5
$ws.Range("L132").Value = 15894
$ws.Range("M132").Value = -1317.5
$ws.Range("N132").Value = -20954

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 4416.3335
$ws.Range("J86").Value = 5470.8
$ws.Range("L86").Value = 5470.8
$ws.Range("N86").Value = -7716.8
$ws.Range("H89").Value = 4416.3335
$ws.Range("J89").Value = 5470.8
$ws.Range("L89").Value = 27354
$ws.Range("N89").Value = -38586
$ws.Range("H94").Value = 662.5
$ws.Range("I94").Value = 665.7
$ws.Range("K94").Value = 665.7
$ws.Range("M94").Value = -214.7
$ws.Range("H105").Value = 1426.125
$ws.Range("J105").Value = 1349
$ws.Range("L105").Value = 1349
$ws.Range("N105").Value = -4843

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1624.4286
$ws.Range("I22").Value = 271.44446
$ws.Range("J22").Value = 4059.8
$ws.Range("K22").Value = 271.44446
$ws.Range("L22").Value = 4059.8
$ws.Range("M22").Value = 78.55554000000001
$ws.Range("N22").Value = -4759.8
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H13").Value = 170.71428
$ws.Range("I13").Value = 35.2
$ws.Range("K13").Value = 105.6
$ws.Range("M13").Value = 62.39999999999999
$ws.Range("H54").Value = 3433.9092
$ws.Range("I54").Value = 3433.9092
$ws.Range("J54").Value = 0
$ws.Range("K54").Value = 10301.7276
$ws.Range("L54").Value = 0
$ws.Range("M54").Value = -9742.7276
$ws.Range("N54").ClearContents()

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 13572232
$ws.Range("I3").Value = 17272818
$ws.Range("J3").Value = 3416.6667
$ws.Range("K3").Value = 17272818
$ws.Range("L3").Value = 3416.6667
$ws.Range("M3").Value = -17272702
$ws.Range("N3").Value = -3648.6667
$ws.Range("H113").Value = 6265.0713
$ws.Range("I113").Value = 3699.8
$ws.Range("K113").Value = 3699.8
$ws.Range("M113").Value = -1529.8

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1436.4286
$ws.Range("J22").Value = 1700
$ws.Range("L22").Value = 1700
$ws.Range("N22").Value = -2290
$ws.Range("H27").Value = 1436.4286
$ws.Range("J27").Value = 1700
$ws.Range("L27").Value = 1700
$ws.Range("N27").Value = -1914
$ws.Range("H43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("N43").ClearContents()
$ws.Range("H46").Value = 6864.8667
$ws.Range("I46").Value = 4994.3335
$ws.Range("J46").Value = 7332.5
$ws.Range("K46").Value = 4994.3335
$ws.Range("L46").Value = 7332.5
$ws.Range("M46").Value = -4806.3335
$ws.Range("N46").Value = -7708.5
$ws.Range("H55").Value = 1052.75
$ws.Range("H93").Value = 1684.1538
$ws.Range("I93").Value = 1754
$ws.Range("J93").Value = 1300
$ws.Range("K93").Value = 1754
$ws.Range("L93").Value = 1300
$ws.Range("M93").Value = -506
$ws.Range("N93").Value = -3796

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H24").Value = 399
$ws.Range("J24").Value = 399
$ws.Range("L24").Value = 399
$ws.Range("N24").Value = -859
$ws.Range("H33").Value = 13504.5
$ws.Range("I33").Value = 11018
$ws.Range("K33").Value = 11018
$ws.Range("M33").Value = -10768
$ws.Range("H36").Value = 13504.5
$ws.Range("I36").Value = 11018
$ws.Range("K36").Value = 11018
$ws.Range("M36").Value = -10768
$ws.Range("H81").Value = 974.5
$ws.Range("I81").Value = 950
$ws.Range("K81").Value = 1900
$ws.Range("M81").Value = -839
$ws.Range("H84").Value = 974.5
$ws.Range("I84").Value = 950
$ws.Range("K84").Value = 9500
$ws.Range("M84").Value = -4196
$ws.Range("H94").Value = 7079583.5
$ws.Range("J94").Value = 7079583.5
$ws.Range("L94").Value = 7079583.5
$ws.Range("N94").Value = -7081385.5
$ws.Range("H107").Value = 503.83334
$ws.Range("I107").Value = 503.83334
$ws.Range("K107").Value = 1511.50002
$ws.Range("M107").Value = 408.4999800000001

